# Update loading_percent values for the 380 kV case (rows 2-25).
# Columns B,D,E,F,G,I,L,N,O change; C,H,J,K,M stay at 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ "B" = 14.97770889801861; "D" = 3.234510661840577; "E" = 18.418478947062; "F" = 21.4860049255458; "G" = 3.583448579622964; "I" = 22.40359339487079; "L" = 9.441462969642364; "N" = 17.76107303827146; "O" = 18.65945141104023 }
    3 = @{ "B" = 14.66425713582991; "D" = 3.225671528957201; "E" = 18.29573303088289; "F" = 21.18169610203237; "G" = 3.586025822461743; "I" = 22.55863817489225; "L" = 9.286277363253578; "N" = 17.75882083674559; "O" = 18.50459139237228 }
    4 = @{ "B" = 14.47015359799598; "D" = 3.220175537588557; "E" = 18.22062115671494; "F" = 21.00015353189779; "G" = 3.587693434609639; "I" = 22.6586139566817; "L" = 9.191104031185096; "N" = 17.75955815525293; "O" = 18.41488877812746 }
    5 = @{ "B" = 14.39074683881289; "D" = 3.217919092348814; "E" = 18.19009584480417; "F" = 20.92760830132575; "G" = 3.588394489439751; "I" = 22.70055981222025; "L" = 9.152395473107223; "N" = 17.7603943845475; "O" = 18.37972527543774 }
    6 = @{ "B" = 14.37754594337595; "D" = 3.217543412125485; "E" = 18.18503273818511; "F" = 20.9156517707063; "G" = 3.588512199066686; "I" = 22.70759776722486; "L" = 9.145973889999345; "N" = 17.76056567298603; "O" = 18.37397146024414 }
    7 = @{ "B" = 14.46908379634478; "D" = 3.220145173690677; "E" = 18.22020911859128; "F" = 20.9991692208755; "G" = 3.587702802175928; "I" = 22.65917476983053; "L" = 9.190581626061244; "N" = 17.75956726026499; "O" = 18.41440887105241 }
    8 = @{ "B" = 14.87003368076345; "D" = 3.231477491279989; "E" = 18.37611497456589; "F" = 21.38003988402463; "G" = 3.584319578044996; "I" = 22.45606336459074; "L" = 9.387960141112933; "N" = 17.75985804914194; "O" = 18.60496213366215 }
    9 = @{ "B" = 15.63883325916875; "D" = 3.253133444058929; "E" = 18.68301148757191; "F" = 22.16444089514406; "G" = 3.578357653563752; "I" = 22.0955083830682; "L" = 9.773852895700202; "N" = 17.77713181660327; "O" = 19.01954607718686 }
    10 = @{ "B" = 16.1874098726353; "D" = 3.268670719216701; "E" = 18.90805136233534; "F" = 22.75737941665145; "G" = 3.574382897600391; "I" = 21.85338950982447; "L" = 10.05398316901023; "N" = 17.79984051230942; "O" = 19.34650207596431 }
    11 = @{ "B" = 16.43235215232298; "D" = 3.275652029511948; "E" = 19.01008826231665; "F" = 23.02940969600635; "G" = 3.572661755212656; "I" = 21.7481409707198; "L" = 10.18014676102371; "N" = 17.81230858345413; "O" = 19.49950812572564 }
    12 = @{ "B" = 16.52435937335301; "D" = 3.278282682416508; "E" = 19.04865819812283; "F" = 23.13264222367547; "G" = 3.572022439296995; "I" = 21.70898596785781; "L" = 10.22769636332203; "N" = 17.81733384409631; "O" = 19.55801028633869 }
    13 = @{ "B" = 16.50457848438903; "D" = 3.27771671204207; "E" = 19.04035488504346; "F" = 23.11040116895531; "G" = 3.572159575046811; "I" = 21.71738760152669; "L" = 10.217466417546; "N" = 17.81623810731585; "O" = 19.54538664137488 }
    14 = @{ "B" = 16.43993707853521; "D" = 3.275868718982589; "E" = 19.01326294167049; "F" = 23.03789889972044; "G" = 3.572608909302129; "I" = 21.74490564760127; "L" = 10.18406346217357; "N" = 17.81271594537531; "O" = 19.5043101709747 }
    15 = @{ "B" = 16.40024264568098; "D" = 3.274735058114351; "E" = 18.99665870619284; "F" = 22.99351466381056; "G" = 3.572885758040941; "I" = 21.76185235487635; "L" = 10.16357252791476; "N" = 17.8105979852799; "O" = 19.47922128154671 }
    16 = @{ "B" = 16.1713022428781; "D" = 3.268212682710304; "E" = 18.90137446841155; "F" = 22.73963840582137; "G" = 3.57449712316859; "I" = 21.86036592516016; "L" = 10.04570879707771; "N" = 17.79906841296505; "O" = 19.3365842937821 }
    17 = @{ "B" = 16.02961259925851; "D" = 3.264188823318869; "E" = 18.84281987596817; "F" = 22.58440346715834; "G" = 3.575507877428389; "I" = 21.92205161826044; "L" = 9.973046742390181; "N" = 17.7925404581161; "O" = 19.25013847282864 }
    18 = @{ "B" = 15.94768743403514; "D" = 3.2618662808196; "E" = 18.8091105579181; "F" = 22.49533935039992; "G" = 3.576097428766398; "I" = 21.95799229749872; "L" = 9.931136472332245; "N" = 17.78898718902777; "O" = 19.20082243685037 }
    19 = @{ "B" = 15.91987796167489; "D" = 3.261078529241065; "E" = 18.79769260802052; "F" = 22.46522544892134; "G" = 3.576298449769554; "I" = 21.97024041783637; "L" = 9.916927635546786; "N" = 17.78781881571366; "O" = 19.18419606281671 }
    20 = @{ "B" = 16.04474073665687; "D" = 3.264618014931203; "E" = 18.84905639771413; "F" = 22.60090622566269; "G" = 3.575399433529127; "I" = 21.9154374134995; "L" = 9.980794168915477; "N" = 17.79321454714423; "O" = 19.25929920489032 }
    21 = @{ "B" = 16.45894474763094; "D" = 3.276411877058107; "E" = 19.02122254807885; "F" = 23.05918945296873; "G" = 3.572476591756045; "I" = 21.73680394897557; "L" = 10.19388118382743; "N" = 17.81374227159118; "O" = 19.51636049183892 }
    22 = @{ "B" = 16.72526115159101; "D" = 3.284043712680996; "E" = 19.13332993017612; "F" = 23.3599404595392; "G" = 3.570638844657076; "I" = 21.62413727006036; "L" = 10.33181098158347; "N" = 17.82892778755051; "O" = 19.68761640864273 }
    23 = @{ "B" = 16.58355118595307; "D" = 3.279977607971079; "E" = 19.07354083485701; "F" = 23.19934656438044; "G" = 3.571613072987081; "I" = 21.6838972734538; "L" = 10.25833119239892; "N" = 17.82066229128463; "O" = 19.59593386221593 }
    24 = @{ "B" = 16.03790275040051; "D" = 3.264424005914648; "E" = 18.84623700643336; "F" = 22.59344475152425; "G" = 3.575448434672867; "I" = 21.91842621015311; "L" = 9.977291977308573; "N" = 17.79290916905228; "O" = 19.25515644144503 }
    25 = @{ "B" = 15.43332073677889; "D" = 3.247338472930773; "E" = 18.59998127413795; "F" = 21.94886467259946; "G" = 3.579898980736789; "I" = 22.1890309788004; "L" = 9.669864498599303; "N" = 17.77713181660327; "O" = 18.90326816083252 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}
